# Refresh crypto price/volume snapshot (coinranking.com) per the
# scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "589.02"); Excel
# auto-converts such literals when assigned straight to .Value, which would
# silently turn them into numbers and drop formatting like "70.327.12" or
# the trailing zeros in "1.00". Force Text (@) just for the assignment, then
# restore the Normal style so no stray formatting is left behind.
function Set-TextValue($cell, $text) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "70.327.12"
$ws.Range("E2").Value = "  +4.40%  "
Set-TextValue "D3" "3.610.53"
$ws.Range("E3").Value = "  +4.34%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue "D5" "589.02"
Set-TextValue "D6" "191.06"
$ws.Range("E6").Value = "  +3.36%  "
Set-TextValue "D7" "0.643"
$ws.Range("E7").Value = "  +1.26%  "
Set-TextValue "D8" "3.605.45"
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("E11").Value = "  +2.22%  "
Set-TextValue "D12" "58.10"
$ws.Range("E12").Value = "  +4.40%  "
Set-TextValue "D13" "0.0000290"
$ws.Range("E13").Value = "  +3.04%  "
Set-TextValue "D14" "9.79"
$ws.Range("E14").Value = "  +4.11%  "
Set-TextValue "D15" "4.190.81"
$ws.Range("E15").Value = "  +4.58%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "3.614.72"
$ws.Range("E16").Value = "  +4.57%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D17" "19.40"
$ws.Range("E17").Value = "  +4.44%  "
Set-TextValue "D18" "70.259.44"
$ws.Range("E18").Value = "  +4.29%  "
Set-TextValue "D19" "12.49"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  +3.72%  "
Set-TextValue "D22" "492.06"
$ws.Range("E22").Value = "  +0.67%  "
Set-TextValue "D23" "17.34"
$ws.Range("E23").Value = "  +14.69%  "
$ws.Range("E24").Value = "  +8.54%  "
Set-TextValue "D25" "4.47"
$ws.Range("E25").Value = "  +6.56%  "
Set-TextValue "D26" "90.80"
$ws.Range("E26").Value = "  +0.56%  "
Set-TextValue "D27" "3.11"
$ws.Range("E27").Value = "  +4.99%  "
Set-TextValue "D28" "11.10"
$ws.Range("E28").Value = "  +0.94%  "
Set-TextValue "D29" "9.43"
$ws.Range("E29").Value = "  +5.17%  "
Set-TextValue "D30" "32.45"
$ws.Range("E30").Value = "  +2.60%  "
Set-TextValue "D31" "7.51"
$ws.Range("E31").Value = "  +7.49%  "
Set-TextValue "D32" "630.37"
$ws.Range("E32").Value = "  +6.12%  "
Set-TextValue "D33" "12.24"
$ws.Range("E33").Value = "  +5.15%  "
$ws.Range("E34").Value = "  +6.88%  "
$ws.Range("E35").Value = "  +2.66%  "
Set-TextValue "D36" "0.0₃0819"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D37" "0.404"
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D38" "1.00"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D39" "37.99"
$ws.Range("E39").Value = "  +3.76%  "
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("E41").Value = "  -0.92%  "
Set-TextValue "D42" "3.303.72"
$ws.Range("E42").Value = "  +5.16%  "
Set-TextValue "D43" "3.09"
$ws.Range("E43").Value = "  +5.69%  "
Set-TextValue "D44" "0.0445"
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("E45").Value = "  +1.84%  "
Set-TextValue "D46" "3.29"
$ws.Range("E46").Value = "  +0.89%  "
Set-TextValue "D47" "0.137"
$ws.Range("E47").Value = "  +1.50%  "
Set-TextValue "D48" "9.14"
$ws.Range("E48").Value = "  +4.06%  "
Set-TextValue "D49" "2.72"
$ws.Range("E49").Value = "  -3.13%  "
Set-TextValue "D50" "3.32"
$ws.Range("E50").Value = "  +6.04%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  -0.12%  "
